$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing the old row 5 (and anything below) down to row 6.
$ws.Rows.Item(5).Insert()

# The old row 5 data is now row 6 and stays unchanged (it already carried over via Insert).
# Populate the newly inserted row 5 with the new weekly record.
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(5, 3).Value = "Coquimbo"
$ws.Cells.Item(5, 4).Value = 44585
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100101
$ws.Cells.Item(5, 8).Value = "Berries"
$ws.Cells.Item(5, 9).Value = 100101004
$ws.Cells.Item(5, 10).Value = "Frambuesa"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 160
$ws.Cells.Item(5, 14).Value = 6500
$ws.Cells.Item(5, 15).Value = 7000
$ws.Cells.Item(5, 16).Value = 6750
$ws.Cells.Item(5, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(5, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(5, 19).Value = 3375
$ws.Cells.Item(5, 20).Value = 2
